$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the text "R40" (a shared string). The edit replaces its
# content with the text "1" - note the leading apostrophe forces Excel to
# store it as text (shared string) instead of coercing it to the number 1.
$ws.Range("B11").Value = "'1"
